# Inserção de ícones no menu
# Fix a typo, add a yellow highlight color, insert new testing-related
# tasks and move one existing task into the "in progress" (yellow) group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 4 blank rows before row 37. This pushes the current rows
#    37-41 ("Todas/CKEditor", "Posts/Tratar dimensionamento",
#    "Postagens/Criar paginação", "General/Tratamento de erros",
#    "Todas/Criar unit tests") down to rows 41-45, and leaves rows
#    37-40 free for the new/re-ordered content.
# ------------------------------------------------------------------
$ws.Range("A37:A40").EntireRow.Insert()

# ------------------------------------------------------------------
# 2. Fix the typo on row 35 ("aterações"/"images" -> "alterações"/"imagens")
# ------------------------------------------------------------------
$ws.Range("B35").Value = "Ao aplicar alterações, não deve alterar imagens que não estejam selecionadas"

# ------------------------------------------------------------------
# 3. Re-colour rows 35-36 green (fill RGB FF00B050) - already red, now
#    marked as belonging to the new "testing" colour group.
# ------------------------------------------------------------------
$ws.Range("D35:D36").Interior.Color = 5287936

# ------------------------------------------------------------------
# 4. Fill in the three brand-new rows (37-39), all green as well.
# ------------------------------------------------------------------
$ws.Range("A37").Value = "Post"
$ws.Range("B37").Value = "Analisar a possibilidade de inserir ícones quando for postar algo"
$ws.Range("C37").Value = "Nova Implementação"

$ws.Range("A38").Value = "Login"
$ws.Range("B38").Value = "Verificar mensagem de erro que acontece após fazer o login"
$ws.Range("C38").Value = "Bug"

$ws.Range("A39").Value = "Postar"
$ws.Range("B39").Value = "Ajustar tamanho dos botões no mobile"
$ws.Range("C39").Value = "Bug"

$ws.Range("D37:D39").Interior.Color = 5287936

# ------------------------------------------------------------------
# 5. Rows 40-45: the "Criar paginação" task is pulled up and
#    highlighted yellow (RGB FFFF00); the remaining four pre-existing
#    rows (CKEditor / Tratar dimensionamento / Tratamento de erros /
#    Criar unit tests) follow in their original relative order
#    (still red), and a brand new row 45 is appended for "Testes de
#    segurança" (also red).
# ------------------------------------------------------------------
$ws.Range("A40").Value = "Postagens"
$ws.Range("B40").Value = "Criar páginação"
$ws.Range("C40").Value = "Nova Implementação"
$ws.Range("D40").Interior.Color = 65535

$ws.Range("A41").Value = "Todas"
$ws.Range("B41").Value = "Atualizar para CKEditor formulários de textos fixos no site"
$ws.Range("C41").Value = "Alteração"
$ws.Range("D41").Interior.Color = 255

$ws.Range("A42").Value = "Posts"
$ws.Range("B42").Value = "Tratar dimensionamento das imagens"
$ws.Range("C42").Value = "Nova Implementação"
$ws.Range("D42").Interior.Color = 255

$ws.Range("A43").Value = "General"
$ws.Range("B43").Value = "Tratamento de erros"
$ws.Range("C43").Value = "Nova Implementação"
$ws.Range("D43").Interior.Color = 255

$ws.Range("A44").Value = "Todas"
$ws.Range("B44").Value = "Criar unit tests"
$ws.Range("C44").Value = "Nova Implementação"
$ws.Range("D44").Interior.Color = 255

$ws.Range("A45").Value = "Todas"
$ws.Range("B45").Value = "Testes de segurança"
$ws.Range("C45").Value = "Testes"
$ws.Range("D45").Interior.Color = 255

# ------------------------------------------------------------------
# 7. Update the view so the selection/scroll matches the edited area.
# ------------------------------------------------------------------
$ws.Range("B40").Select()
